$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix up the header/migration rows: the raw inlineStr values previously carried
# a bunch of concatenated leftover data from the subentity migration; now that
# migration is complete so they just hold the plain field definitions.
$ws.Range("B2").Value = "name=Type,dataType=text,updateCriteria=true"
$ws.Range("B3").Value = "Type"

# Append the migrated "Type" subentity rows (id, label, visibilityGroups).
$rows = @(
    @(2, "Fighting Weapon", -1),
    @(3, "Ranged Weapon", -1),
    @(4, "Shield", -1),
    @(5, "Armor", -1),
    @(6, "Accessory", -1),
    @(7, "Consumable", -1)
)

$r = 4
foreach ($row in $rows) {
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]

    $ws.Cells.Item($r, 1).Value = $row[0]
    # Column A carries the same bold/bordered/centered style used by the
    # existing id cells (A2/A3) -- copy it across instead of reinventing it.
    $ws.Range("A3").Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)

    $r++
}

$excel.CutCopyMode = $false

Write-Output "done"
